$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values in column D stay as text, matching the
# original inline-string cell type (avoids Excel auto-converting "305.89" etc.
# into a number or "41.926.16" into something else).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.957.21"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.272.09"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.89"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.02"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.76"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.112"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.67"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.624.34"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.36"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.272.25"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.785"
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.875.37"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.75"
$ws.Range("E19").Value = "  +3.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0918"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.19"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.13"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.94"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.00"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.09"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.04"
$ws.Range("E30").Value = "  +2.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.17"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.37"
$ws.Range("E32").Value = "  +3.60%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.25"
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.96"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.89"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.013.45"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.24"
$ws.Range("E44").Value = "  +9.85%  "
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.32"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.92"
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.44"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.56"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.51"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("E51").Value = "  -0.26%  "
